$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Collin Sexton", "PG,SG", "Utah Jazz"),
    @("Buddy Hield", "SG,SF", "Golden State Warriors"),
    @("Gradey Dick", "SG,SF", "Toronto Raptors"),
    @("Anfernee Simons", "PG,SG", "Portland Trail Blazers"),
    @("Andrew Wiggins", "SF,PF", "Milwaukee Bucks"),
    @("Giannis Antetokounmpo", "PF,C", "Minnesota Timberwolves"),
    @("Anthony Edwards", "SG,SF", "Philadelphia 76ers"),
    @("Paul George", "SG,SF,PF", "Memphis Grizzlies"),
    @("Jaren Jackson Jr.", "PF,C", "Denver Nuggets"),
    @("Aaron Gordon", "PF,C", "LA Clippers"),
    @("James Harden", "PG,SG", "LA Clippers"),
    @("Ivica Zubac", "C", "LA Clippers"),
    @("Tyus Jones", "PG", "Phoenix Suns"),
    @("Fred VanVleet", "PG", "Houston Rockets"),
    @("Jaden Ivey", "PG,SG", "Detroit Pistons"),
    @("Khris Middleton", "SF", "Milwaukee Bucks"),
    @("Zion Williamson", "PF,C", "New Orleans Pelicans"),
    @("Jayson Tatum", "SF,PF", "Boston Celtics")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
